$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.281.55'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.897.89'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.694'
$ws.Range('E5').Value = '  +9.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '245.25'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.67'
$ws.Range('E8').Value = '  -4.58%  '
$ws.Range('E9').Value = '  +2.68%  '
$ws.Range('E10').Value = '  +10.99%  '
$ws.Range('E11').Value = '  +1.74%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '2.173.59'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.51'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.704'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '1.899.82'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '35.264.58'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.12'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.57'
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.58'
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  +6.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.62'
$ws.Range('E27').Value = '  -2.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.58'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.131'
$ws.Range('E29').Value = '  +3.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.28'
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.15'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0566'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('E35').Value = '  +16.06%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.10'
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.909'
$ws.Range('E37').Value = '  -4.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.82'
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('E40').Value = '  +10.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0209'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.04'
$ws.Range('E43').Value = '  +5.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '89.41'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').Value = '1.351.58'
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.67'
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.84'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.49'
$ws.Range('E51').Value = '  -2.81%  '
